# Updated cryptos list on Mon Jun  3 15:55:01 UTC 2024 with GitHub Actions
# Refreshes prices / 1h volume % for existing rows and re-ranks two swapped
# coin pairs (Fetch.AI <-> InternetComputer, Arweave <-> OKB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces text storage so decimal-look-alike prices (trailing
# zeros, etc.) keep their exact original formatting instead of being
# auto-coerced to numbers by Excel.
$ws.Range("D2").Value = "69.242.95"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.786.88"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'627.42"
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("D6").Value = "'165.00"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "3.784.89"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "4.423.11"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "3.772.07"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "69.232.87"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "'17.94"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "'7.13"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'469.68"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Value = "'9.65"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "'83.18"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'2.17"
$ws.Range("E26").Value = "  +2.54%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "3.934.44"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "'2.25"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D34").Value = "'28.94"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'9.04"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "3.735.96"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "'0.969"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'0.301"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  +3.66%  "
$ws.Range("D47").Value = "'153.05"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'43.07"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'46.83"
$ws.Range("E49").Value = "  -1.75%  "
$ws.Range("D50").Value = "'1.40"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("E51").Value = "  +0.85%  "
